$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracking_Main")
$v = $ws.Range("A1").Value
Write-Host "Type: $($v.GetType())"
Write-Host ("A1 value: " + $v)
